$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.96264576182254
$ws.Range("C2").Value = 7.892458382824905
$ws.Range("E2").Value = 11.99529766303683
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 45.82342251323082
$ws.Range("H2").Value = 18.62657270360087
$ws.Range("K2").Value = 11.33732425459352
$ws.Range("L2").Value = 10.01407828799498
$ws.Range("M2").Value = 15.37425287799839
$ws.Range("N2").Value = 22.13714351742483
$ws.Range("B3").Value = 14.8075288436562
$ws.Range("C3").Value = 7.85420619193022
$ws.Range("E3").Value = 12.01620645695606
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 45.90474952519263
$ws.Range("H3").Value = 18.67759482141792
$ws.Range("K3").Value = 11.22762221008536
$ws.Range("L3").Value = 10.022624728623
$ws.Range("M3").Value = 15.35998060302919
$ws.Range("N3").Value = 22.2040121925503
$ws.Range("B4").Value = 14.71493774363437
$ws.Range("C4").Value = 7.830169612476418
$ws.Range("E4").Value = 12.03045534371482
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 45.96591837041424
$ws.Range("H4").Value = 18.71172596283878
$ws.Range("K4").Value = 11.16220507417197
$ws.Range("L4").Value = 10.02921077590826
$ws.Range("M4").Value = 15.35370321699419
$ws.Range("N4").Value = 22.24702056860994
$ws.Range("B5").Value = 14.67791436038873
$ws.Range("C5").Value = 7.820238052235501
$ws.Range("E5").Value = 12.03661699158595
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 45.99366175679815
$ws.Range("H5").Value = 18.72633926518057
$ws.Range("K5").Value = 11.13606297237456
$ws.Range("L5").Value = 10.03223167179247
$ws.Range("M5").Value = 15.35177311097317
$ws.Range("N5").Value = 22.26503868610356
$ws.Range("B6").Value = 14.67181061863466
$ws.Range("C6").Value = 7.818580733176163
$ws.Range("E6").Value = 12.03766158637515
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 45.99843838191687
$ws.Range("H6").Value = 18.72880833795676
$ws.Range("K6").Value = 11.13175405467559
$ws.Range("L6").Value = 10.03275365696837
$ws.Range("M6").Value = 15.35149062474323
$ws.Range("N6").Value = 22.26806032666243
$ws.Range("B7").Value = 14.71443551184476
$ws.Range("C7").Value = 7.830036222197495
$ws.Range("E7").Value = 12.03053700356123
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 45.96628113476896
$ws.Range("H7").Value = 18.71192019058354
$ws.Range("K7").Value = 11.16185038749588
$ws.Range("L7").Value = 10.02925015159794
$ws.Range("M7").Value = 15.35367464067487
$ws.Range("N7").Value = 22.24726157399562
$ws.Range("B8").Value = 14.90863650006781
$ws.Range("C8").Value = 7.879384129599195
$ws.Range("E8").Value = 12.00221441939398
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 45.84912781726491
$ws.Range("H8").Value = 18.64358309232848
$ws.Range("K8").Value = 11.29911335073662
$ws.Range("L8").Value = 10.01674758245562
$ws.Range("M8").Value = 15.368817690697
$ws.Range("N8").Value = 22.1597957624118
$ws.Range("B9").Value = 15.30861586750744
$ws.Range("C9").Value = 7.971720553096572
$ws.Range("E9").Value = 11.95785225617213
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 45.70885721190077
$ws.Range("H9").Value = 18.53183277409187
$ws.Range("K9").Value = 11.58240886615556
$ws.Range("L9").Value = 10.0028291344537
$ws.Range("M9").Value = 15.41809760451979
$ws.Range("N9").Value = 22.0036921912826
$ws.Range("B10").Value = 15.6115686537961
$ws.Range("C10").Value = 8.036760908424442
$ws.Range("E10").Value = 11.93205242686741
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 45.6607418151185
$ws.Range("H10").Value = 18.46331976626744
$ws.Range("K10").Value = 11.7974045524056
$ws.Range("L10").Value = 9.999033826294855
$ws.Range("M10").Value = 15.4660371943395
$ws.Range("N10").Value = 21.89831344412104
$ws.Range("B11").Value = 15.75080333677452
$ws.Range("C11").Value = 8.065721491991715
$ws.Range("E11").Value = 11.92178555658458
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 45.65084773933463
$ws.Range("H11").Value = 18.4351066180127
$ws.Range("K11").Value = 11.89632018722383
$ws.Range("L11").Value = 9.998695739522651
$ws.Range("M11").Value = 15.49034331093869
$ws.Range("N11").Value = 21.85237691281158
$ws.Range("B12").Value = 15.80368515447752
$ws.Range("C12").Value = 8.076596229154772
$ws.Range("E12").Value = 11.91810864805767
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 45.64882922505483
$ws.Range("H12").Value = 18.42484814132583
$ws.Range("K12").Value = 11.93390493574449
$ws.Range("L12").Value = 9.99876657031677
$ws.Range("M12").Value = 15.49990175752068
$ws.Range("N12").Value = 21.83526828738458
$ws.Range("B13").Value = 15.79228989338956
$ws.Range("C13").Value = 8.074258282166705
$ws.Range("E13").Value = 11.91889116133905
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 45.64918704514159
$ws.Range("H13").Value = 18.42703856809299
$ws.Range("K13").Value = 11.92580523082825
$ws.Range("L13").Value = 9.998742484858113
$ws.Range("M13").Value = 15.49782750264106
$ws.Range("N13").Value = 21.83894021041237
$ws.Range("B14").Value = 15.75515108669902
$ws.Range("C14").Value = 8.066618017986853
$ws.Range("E14").Value = 11.92147883004474
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 45.65064702438917
$ws.Range("H14").Value = 18.4342541220117
$ws.Range("K14").Value = 11.89940993651974
$ws.Range("L14").Value = 9.998697586238803
$ws.Range("M14").Value = 15.49112262016402
$ws.Range("N14").Value = 21.85096363981952
$ws.Range("B15").Value = 15.73242149400155
$ws.Range("C15").Value = 8.06192609721136
$ws.Range("E15").Value = 11.92309130939119
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 45.6517664439275
$ws.Range("H15").Value = 18.43872924869184
$ws.Range("K15").Value = 11.8832577046442
$ws.Range("L15").Value = 9.998695956737915
$ws.Range("M15").Value = 15.48706166598218
$ws.Range("N15").Value = 21.85836561927105
$ws.Range("B16").Value = 15.60249384236804
$ws.Range("C16").Value = 8.034855494473273
$ws.Range("E16").Value = 11.93275292990493
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 45.66163008479305
$ws.Range("H16").Value = 18.46522304077382
$ws.Range("K16").Value = 11.79095979800996
$ws.Range("L16").Value = 9.999083793175476
$ws.Range("M16").Value = 15.46449858329707
$ws.Range("N16").Value = 21.90135566521829
$ws.Range("B17").Value = 15.52311720419057
$ws.Range("C17").Value = 8.018086711976963
$ws.Range("E17").Value = 11.93905614759723
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 45.67075584427413
$ws.Range("H17").Value = 18.48223300249341
$ws.Range("K17").Value = 11.73459992833081
$ws.Range("L17").Value = 9.999676781794077
$ws.Range("M17").Value = 15.45129326054816
$ws.Range("N17").Value = 21.92824027991178
$ws.Range("B18").Value = 15.47759834331697
$ws.Range("C18").Value = 8.008382886495582
$ws.Range("E18").Value = 11.94281994252124
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 45.67713352342005
$ws.Range("H18").Value = 18.49229467680406
$ws.Range("K18").Value = 11.70228992730211
$ws.Range("L18").Value = 10.00014858975756
$ws.Range("M18").Value = 15.44393336610713
$ws.Range("N18").Value = 21.94389199191644
$ws.Range("B19").Value = 15.46221133353774
$ws.Range("C19").Value = 8.005087285888756
$ws.Range("E19").Value = 11.94411807235699
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 45.67948664178802
$ws.Range("H19").Value = 18.49574911765391
$ws.Range("K19").Value = 11.69136963535823
$ws.Range("L19").Value = 10.00033081304449
$ws.Range("M19").Value = 15.44148201884819
$ws.Range("N19").Value = 21.94922378516851
$ws.Range("B20").Value = 15.53155320283916
$ws.Range("C20").Value = 8.019877883665607
$ws.Range("E20").Value = 11.9383708433151
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 45.66966753801627
$ws.Range("H20").Value = 18.48039348809927
$ws.Range("K20").Value = 11.74058872498068
$ws.Range("L20").Value = 9.999600131936106
$ws.Range("M20").Value = 15.4526746557299
$ws.Range("N20").Value = 21.92535887733739
$ws.Range("B21").Value = 15.76605577934936
$ws.Range("C21").Value = 8.068864660942422
$ws.Range("E21").Value = 11.92071304762772
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 45.65017126954544
$ws.Range("H21").Value = 18.43212319452939
$ws.Range("K21").Value = 11.90715967099558
$ws.Range("L21").Value = 9.998705383829549
$ws.Range("M21").Value = 15.49308243385572
$ws.Range("N21").Value = 21.84742429972882
$ws.Range("B22").Value = 15.92020823086661
$ws.Range("C22").Value = 8.100343188280187
$ws.Range("E22").Value = 11.91040193432072
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 45.6475031432892
$ws.Range("H22").Value = 18.40305454602404
$ws.Range("K22").Value = 12.01675191721163
$ws.Range("L22").Value = 9.999279284721533
$ws.Range("M22").Value = 15.52155353614052
$ws.Range("N22").Value = 21.79815920166839
$ws.Range("B23").Value = 15.83786818006271
$ws.Range("C23").Value = 8.083592281998895
$ws.Range("E23").Value = 11.91579282665886
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 45.64800457417442
$ws.Range("H23").Value = 18.41834206747328
$ws.Range("K23").Value = 11.95820450196571
$ws.Range("L23").Value = 9.998867246019204
$ws.Range("M23").Value = 15.50617100558289
$ws.Range("N23").Value = 21.82430052091471
$ws.Range("B24").Value = 15.52773892342623
$ws.Range("C24").Value = 8.019068291771074
$ws.Range("E24").Value = 11.93868023330046
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 45.67015603799971
$ws.Range("H24").Value = 18.48122425283951
$ws.Range("K24").Value = 11.73788090083235
$ws.Range("L24").Value = 9.999634377567215
$ws.Range("M24").Value = 15.45204940386275
$ws.Range("N24").Value = 21.92666095069559
$ws.Range("B25").Value = 15.19862530120621
$ws.Range("C25").Value = 7.947225752049714
$ws.Range("E25").Value = 11.96865883064075
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 45.73718077336151
$ws.Range("H25").Value = 18.55967909661729
$ws.Range("K25").Value = 11.50443683217814
$ws.Range("L25").Value = 10.00546266649849
$ws.Range("M25").Value = 15.40268928638534
$ws.Range("N25").Value = 22.04428101818493
